$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-7) were re-associated with a different
# combination of date / variety / quality / volume / price / region values.
# Target state per row (columns D, K, L, M, N, O, P, Q, R, S, T):
$rows = @{
    2 = @{ D = 44160; K = "Castle Brite";  L = "Primera"; M = 240; N = 20500; O = 21000; P = 20750; Q = "`$/caja 15 kilos"; R = "Región Metropolitana";  S = 1383; T = 15 }
    3 = @{ D = 44175; K = "Castle Brite";  L = "Primera"; M = 300; N = 21000; O = 22000; P = 21500; Q = "`$/caja 18 kilos"; R = "Región Metropolitana";  S = 1194; T = 18 }
    4 = @{ D = 44189; K = "Dina";          L = "Especial"; M = 120; N = 23500; O = 24000; P = 23750; Q = "`$/caja 18 kilos"; R = "Región de O'Higgins";  S = 1319; T = 18 }
    5 = @{ D = 44189; K = "Dina";          L = "Primera"; M = 200; N = 21500; O = 22000; P = 21750; Q = "`$/caja 18 kilos"; R = "Región de O'Higgins";  S = 1208; T = 18 }
    6 = @{ D = 44161; K = "Dina";          L = "Primera"; M = 300; N = 20000; O = 20500; P = 20250; Q = "`$/caja 15 kilos"; R = "Región Metropolitana";  S = 1350; T = 15 }
    7 = @{ D = 44161; K = "Dina";          L = "Segunda"; M = 100; N = 18000; O = 18500; P = 18250; Q = "`$/caja 15 kilos"; R = "Región Metropolitana";  S = 1217; T = 15 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
